$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2005314.6
$ws.Range("J17").Value = 2005314.6
$ws.Range("L17").Value = 6015943.800000001
$ws.Range("N17").Value = -6016279.800000001

$ws.Range("H74").Value = 4854.2856
$ws.Range("I74").Value = 4745
$ws.Range("K74").Value = 4745
$ws.Range("M74").Value = -3809

$ws.Range("H76").Value = 3658.8125
$ws.Range("I76").Value = 3003
$ws.Range("K76").Value = 3003
$ws.Range("M76").Value = -2688

$ws.Range("H77").Value = 4854.2856
$ws.Range("I77").Value = 4745
$ws.Range("K77").Value = 23725
$ws.Range("M77").Value = -19045

$ws.Range("H79").Value = 3658.8125
$ws.Range("I79").Value = 3003
$ws.Range("K79").Value = 3003
$ws.Range("M79").Value = -1911

$ws.Range("H107").Value = 376.25925
$ws.Range("I107").Value = 209.22728
$ws.Range("K107").Value = 209.22728
$ws.Range("M107").Value = 1710.77272

$ws.Range("H138").Value = 1624.51
$ws.Range("I138").Value = 1132.3962
$ws.Range("J138").Value = 2179.4468
$ws.Range("K138").Value = 3397.188599999999
$ws.Range("L138").Value = 6538.340400000001
$ws.Range("M138").Value = 1742.811400000001
$ws.Range("N138").Value = -16818.3404

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 28240.334
$ws.Range("I110").Value = 29848.588
$ws.Range("K110").Value = 29848.588
$ws.Range("M110").Value = -27803.588

$ws.Range("H122").Value = 1874.8125
$ws.Range("I122").Value = 1275.8889
$ws.Range("J122").Value = 2644.8572
$ws.Range("K122").Value = 3827.6667
$ws.Range("L122").Value = 7934.571599999999
$ws.Range("M122").Value = -1377.6667
$ws.Range("N122").Value = -12834.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 46776.89
$ws.Range("I133").Value = 19999
$ws.Range("J133").Value = 50124.125
$ws.Range("K133").Value = 19999
$ws.Range("L133").Value = 50124.125
$ws.Range("M133").Value = -14939
$ws.Range("N133").Value = -60244.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1525.5555
$ws.Range("I16").Value = 1591.25
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1591.25
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1304.25
$ws.Range("N16").Value = -1574

$ws.Range("H99").Value = 7963.5557
$ws.Range("I99").Value = 2087
$ws.Range("J99").Value = 19716.666
$ws.Range("K99").Value = 2087
$ws.Range("L99").Value = 19716.666
$ws.Range("M99").Value = -589
$ws.Range("N99").Value = -22712.666

$ws.Range("H113").Value = 1525.5555
$ws.Range("I113").Value = 1591.25
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1591.25
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 578.75
$ws.Range("N113").Value = -5340

$ws.Range("H122").Value = 1115612
$ws.Range("I122").Value = 3926.6667
$ws.Range("J122").Value = 1671454.6
$ws.Range("K122").Value = 11780.0001
$ws.Range("L122").Value = 5014363.800000001
$ws.Range("M122").Value = -9330.000100000001
$ws.Range("N122").Value = -5019263.800000001

$ws.Range("H126").Value = 7963.5557
$ws.Range("I126").Value = 2087
$ws.Range("J126").Value = 19716.666
$ws.Range("K126").Value = 6261
$ws.Range("L126").Value = 59149.99800000001
$ws.Range("M126").Value = -3791
$ws.Range("N126").Value = -64089.99800000001

$ws.Range("H133").Value = 32825
$ws.Range("J133").Value = 32825
$ws.Range("L133").Value = 32825
$ws.Range("N133").Value = -37885

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 732134.8
$ws.Range("I5").Value = 462.3846
$ws.Range("K5").Value = 1387.1538
$ws.Range("M5").Value = -1275.1538

$ws.Range("H8").Value = 65.333336
$ws.Range("I8").Value = 65.333336
$ws.Range("K8").Value = 196.000008
$ws.Range("M8").Value = -57.00000800000001

$ws.Range("H33").Value = 72.2
$ws.Range("I33").Value = 78.5
$ws.Range("J33").Value = 68
$ws.Range("K33").Value = 471
$ws.Range("L33").Value = 408
$ws.Range("M33").Value = -188
$ws.Range("N33").Value = -974

$ws.Range("H131").Value = 871.5816
$ws.Range("J131").Value = 883.48956
$ws.Range("L131").Value = 2650.46868
$ws.Range("N131").Value = -12730.46868

$ws.Range("H135").Value = 732134.8
$ws.Range("I135").Value = 462.3846
$ws.Range("K135").Value = 4161.4614
$ws.Range("M135").Value = -1626.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5358.8164
$ws.Range("J70").Value = 5760.517
$ws.Range("L70").Value = 5760.517
$ws.Range("N70").Value = -6300.517

$ws.Range("H73").Value = 5358.8164
$ws.Range("J73").Value = 5760.517
$ws.Range("L73").Value = 5760.517
$ws.Range("N73").Value = -7632.517

$ws.Range("H80").Value = 2335
$ws.Range("I80").Value = 2411.6667
$ws.Range("J80").Value = 2296.6667
$ws.Range("K80").Value = 2411.6667
$ws.Range("L80").Value = 2296.6667
$ws.Range("M80").Value = -1413.6667
$ws.Range("N80").Value = -4292.6667

$ws.Range("H83").Value = 2335
$ws.Range("I83").Value = 2411.6667
$ws.Range("J83").Value = 2296.6667
$ws.Range("K83").Value = 12058.3335
$ws.Range("L83").Value = 11483.3335
$ws.Range("M83").Value = -7066.333500000001
$ws.Range("N83").Value = -21467.3335

$ws.Range("H102").Value = 7696190
$ws.Range("I102").Value = 10992698
$ws.Range("J102").Value = 4338
$ws.Range("K102").Value = 10992698
$ws.Range("L102").Value = 4338
$ws.Range("M102").Value = -10991076
$ws.Range("N102").Value = -7582

$ws.Range("H107").Value = 1167.8334
$ws.Range("I107").Value = 1861.6666
$ws.Range("J107").Value = 474
$ws.Range("K107").Value = 1861.6666
$ws.Range("L107").Value = 474
$ws.Range("M107").Value = 58.33339999999998
$ws.Range("N107").Value = -4314

$ws.Range("H122").Value = 92283.39
$ws.Range("I122").Value = 134099.42
$ws.Range("J122").Value = 4005.111
$ws.Range("K122").Value = 402298.26
$ws.Range("L122").Value = 12015.333
$ws.Range("M122").Value = -399848.26
$ws.Range("N122").Value = -16915.333

$ws.Range("H126").Value = 2237.0386
$ws.Range("I126").Value = 2116.9375
$ws.Range("J126").Value = 2429.2
$ws.Range("K126").Value = 6350.8125
$ws.Range("L126").Value = 7287.599999999999
$ws.Range("M126").Value = -3880.8125
$ws.Range("N126").Value = -12227.6

$ws.Range("H132").Value = 2498.3572
$ws.Range("I132").Value = 2070.6316
$ws.Range("J132").Value = 3401.3333
$ws.Range("K132").Value = 6211.8948
$ws.Range("L132").Value = 10203.9999
$ws.Range("M132").Value = -3681.8948
$ws.Range("N132").Value = -15263.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 649.65515
$ws.Range("I22").Value = 555.7143
$ws.Range("J22").Value = 737.3333
$ws.Range("K22").Value = 555.7143
$ws.Range("L22").Value = 737.3333
$ws.Range("M22").Value = -260.7143
$ws.Range("N22").Value = -1327.3333

$ws.Range("H27").Value = 649.65515
$ws.Range("I27").Value = 555.7143
$ws.Range("J27").Value = 737.3333
$ws.Range("K27").Value = 555.7143
$ws.Range("L27").Value = 737.3333
$ws.Range("M27").Value = -448.7143
$ws.Range("N27").Value = -951.3333

$ws.Range("H61").Value = 168884.33
$ws.Range("I61").Value = 202501.2
$ws.Range("J61").Value = 800
$ws.Range("K61").Value = 202501.2
$ws.Range("L61").Value = 800
$ws.Range("M61").Value = -202299.2
$ws.Range("N61").Value = -1204

$ws.Range("H113").Value = 168884.33
$ws.Range("I113").Value = 202501.2
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 202501.2
$ws.Range("L113").Value = 800
$ws.Range("M113").Value = -200331.2
$ws.Range("N113").Value = -5140

$ws.Range("H122").Value = 27779744
$ws.Range("I122").Value = 55556810
$ws.Range("K122").Value = 166670430
$ws.Range("M122").Value = -166667980

$ws.Range("H133").Value = 42577.75
$ws.Range("J133").Value = 42577.75
$ws.Range("L133").Value = 42577.75
$ws.Range("N133").Value = -47637.75

$ws.Range("H136").Value = 23813154
$ws.Range("I136").Value = 3775.8333
$ws.Range("J136").Value = 166669420
$ws.Range("K136").Value = 11327.4999
$ws.Range("L136").Value = 500008260
$ws.Range("M136").Value = -8777.499899999999
$ws.Range("N136").Value = -500013360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 48730
$ws.Range("J64").Value = 48730
$ws.Range("L64").Value = 48730
$ws.Range("N64").Value = -49226

$ws.Range("H67").Value = 48730
$ws.Range("J67").Value = 48730
$ws.Range("L67").Value = 48730
$ws.Range("N67").Value = -50446

$ws.Range("H108").Value = 23353.75
$ws.Range("J108").Value = 23353.75
$ws.Range("L108").Value = 23353.75
$ws.Range("N108").Value = -31033.75

$ws.Range("H113").Value = 1336.7
$ws.Range("I113").Value = 1711
$ws.Range("J113").Value = 463.33334
$ws.Range("K113").Value = 5133
$ws.Range("L113").Value = 1390.00002
$ws.Range("M113").Value = -2963
$ws.Range("N113").Value = -5730.000019999999

$ws.Range("H122").Value = 102460.5
$ws.Range("I122").Value = 202600
$ws.Range("K122").Value = 607800
$ws.Range("M122").Value = -605350

$ws.Range("H132").Value = 1946.7609
$ws.Range("I132").Value = 1233.4828
$ws.Range("J132").Value = 3163.5293
$ws.Range("K132").Value = 3700.4484
$ws.Range("L132").Value = 9490.5879
$ws.Range("M132").Value = -1170.4484
$ws.Range("N132").Value = -14550.5879

$ws.Range("H136").Value = 2565.7812
$ws.Range("I136").Value = 900
$ws.Range("J136").Value = 3705.5264
$ws.Range("K136").Value = 2700
$ws.Range("L136").Value = 11116.5792
$ws.Range("M136").Value = -150
$ws.Range("N136").Value = -16216.5792

